$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.085.44"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "1.637.01"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").Value = "'216.61"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").Value = "'0.516"
$ws.Range("E6").Value = "  +1.74%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").Value = "'0.0624"
$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("D10").Value = "'19.89"
$ws.Range("E10").Value = "  +0.13%  "

$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  +0.03%  "

$ws.Range("D12").Value = "1.868.05"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").Value = "1.631.48"
$ws.Range("E13").Value = "  -0.48%  "

$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("D15").Value = "'0.540"
$ws.Range("E15").Value = "  +2.05%  "

$ws.Range("D16").Value = "'66.59"
$ws.Range("E16").Value = "  -0.87%  "

$ws.Range("D17").Value = "27.098.79"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").Value = "0.0₃0738"
$ws.Range("E18").Value = "  +1.21%  "

$ws.Range("D19").Value = "'216.70"
$ws.Range("E19").Value = "  -0.73%  "

$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").Value = "'6.87"
$ws.Range("E21").Value = "  +2.13%  "

$ws.Range("E22").Value = "  +0.42%  "

$ws.Range("D23").Value = "'2.50"
$ws.Range("E23").Value = "  +2.72%  "

$ws.Range("E24").Value = "  -0.67%  "

$ws.Range("D25").Value = "'146.79"
$ws.Range("E25").Value = "  -0.25%  "

$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("D27").Value = "'7.39"
$ws.Range("E27").Value = "  +1.98%  "

$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("E29").Value = "  -0.88%  "

$ws.Range("D30").Value = "'0.0506"
$ws.Range("E30").Value = "  +0.88%  "

$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").Value = "'3.37"
$ws.Range("E32").Value = "  +1.48%  "

$ws.Range("E33").Value = "  +0.49%  "

$ws.Range("D34").Value = "1.301.32"
$ws.Range("E34").Value = "  +2.84%  "

$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("E36").Value = "  +1.42%  "

$ws.Range("E37").Value = "  -0.53%  "

$ws.Range("D38").Value = "'0.854"
$ws.Range("E38").Value = "  +2.23%  "

$ws.Range("E39").Value = "  +1.62%  "

$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("D41").Value = "'0.807"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("E42").Value = "  +5.56%  "

$ws.Range("E43").Value = "  -1.49%  "

$ws.Range("D44").Value = "1.777.18"
$ws.Range("E44").Value = "  -0.12%  "

$ws.Range("D45").Value = "'61.60"
$ws.Range("E45").Value = "  -0.44%  "

$ws.Range("D46").Value = "'91.17"
$ws.Range("E46").Value = "  -0.68%  "

$ws.Range("E47").Value = "  +0.79%  "

$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +1.90%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").Value = "'7.63"
$ws.Range("E50").Value = "  +0.22%  "

$ws.Range("D51").Value = "'0.0956"
$ws.Range("E51").Value = "  -0.30%  "
